$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target dataset: date label, then the three price-index columns.
$data = @(
    @("2018-10", 100.5, 101.9, 108.1),
    @("2018-11", 100.6, 103.1, 111.6),
    @("2018-12", 101.1, 103.3, 111.7),
    @("2018-01", 99.9308, 103.149, 111.7546),
    @("2018-02", 102.4, 104.3, 111.6),
    @("2018-03", 101.1, 105.3, 112),
    @("2018-04", 102.3, 103.9, 111),
    @("2018-05", 103.5, 103.2, 110.3),
    @("2018-06", 102, 105.9, 110.7),
    @("2018-07", 101.9, 106.8, 108.7),
    @("2018-08", 102.6, 105, 114.2),
    @("2018-09", 100.9, 102.8, 112.7),
    @("2019-10", 102.5, 98.3, 106.7),
    @("2019-11", 102.4, 96.2, 108),
    @("2019-12", 102.2, 95.3, 108.6),
    @("2019-01", 102.5, 102.3, 113.5),
    @("2019-02", 102.4, 101, 111.5),
    @("2019-03", 102.4, 100.5, 111.7),
    @("2019-04", 103, 101.1, 110.6),
    @("2019-05", 103.1, 103.1, 112.2),
    @("2019-06", 103.1, 102.5, 113.4),
    @("2019-07", 102.9, 100.2, 112.5),
    @("2019-08", 102.9, 100.3, 107.8),
    @("2019-09", 102.7, 99.3, 109.2),
    @("2020-10", 100, 93.40000000000001, 111),
    @("2020-11", 100, 96.3, 109.7),
    @("2020-12", 100, 100.7, 110.4),
    @("2020-01", 100.6, 94.7, 110.2),
    @("2020-02", 100.6, 94.59999999999999, 114),
    @("2020-03", 100.6, 95.09999999999999, 113.3),
    @("2020-04", 100, 93, 117.3),
    @("2020-05", 100, 89.5, 115.8),
    @("2020-06", 100, 89.2, 112.9),
    @("2020-07", 100, 91.2, 113.3),
    @("2020-08", 100, 90.7, 111.7),
    @("2020-09", 100, 91.2, 112),
    @("2021-10", $null, 206.1, 125.5),
    @("2021-11", $null, 190.5, 130.2),
    @("2021-12", $null, 167.8, 130.8),
    @("2021-01", $null, 109.2, 102.9),
    @("2021-02", $null, 111.1, 103.1),
    @("2021-03", $null, 107, 102),
    @("2021-04", $null, 113.6, 101.3),
    @("2021-05", $null, 130.5, 104.7),
    @("2021-06", $null, 138.3, 106.6),
    @("2021-07", $null, 146.9, 107.8),
    @("2021-08", $null, 158.4, 113.8),
    @("2021-09", $null, 176.8, 116),
    @("2022-10", $null, 83, 107.2),
    @("2022-11", $null, 88.2, 101.7),
    @("2022-12", $null, 97.09999999999999, 104.2),
    @("2022-01", $null, 151.8, 130.2),
    @("2022-02", $null, 145.6, 133),
    @("2022-03", $null, 154.3, 136.7),
    @("2022-04", $null, 154, 131.6),
    @("2022-05", $null, 137.5, 124.6),
    @("2022-06", $null, 131.5, 126.4),
    @("2022-07", $null, 120.6, 124.5),
    @("2022-08", $null, 108.2, 119.8),
    @("2022-09", $null, 96.90000000000001, 115.7),
    @("2023-01", $null, 100.4, 101.5),
    @("2023-02", $null, 100.6, 98.7),
    @("2023-03", $null, 96.8, 97.09999999999999),
    @("2023-04", $null, 90.5, 97.8),
    @("2023-05", $null, 86.59999999999999, 100),
    @("2023-06", $null, 80.2, 98.8),
    @("2023-07", $null, 80.40000000000001, 98.3)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -eq $null) {
        $ws.Cells.Item($r, 2).Value = ""
    } else {
        $ws.Cells.Item($r, 2).Value = $row[1]
    }
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# New rows (50-68) were appended beyond the former A1:D49 extent; copy the
# date-column formatting (bold/centered/bordered style) from an existing
# labeled cell so the newly written dates match the rest of column A.
$ws.Range("A2").Copy()
$ws.Range("A50:A68").PasteSpecial(-4122)
$excel.CutCopyMode = 0